# Auto-generated script applying market-data refresh to Tonberry_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1387.8667
$ws.Range("I38").Value = 1021.8333
$ws.Range("J38").Value = 1631.8889
$ws.Range("K38").Value = 3065.4999
$ws.Range("L38").Value = 4895.6667
$ws.Range("M38").Value = -2693.4999
$ws.Range("N38").Value = -5639.6667
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3350
$ws.Range("H62").Value = 3695.6
$ws.Range("I62").Value = 1917
$ws.Range("J62").Value = 5728.2856
$ws.Range("K62").Value = 1917
$ws.Range("L62").Value = 5728.2856
$ws.Range("M62").Value = -1293
$ws.Range("N62").Value = -6976.2856
$ws.Range("H65").Value = 3695.6
$ws.Range("I65").Value = 1917
$ws.Range("J65").Value = 5728.2856
$ws.Range("K65").Value = 9585
$ws.Range("L65").Value = 28641.428
$ws.Range("M65").Value = -6465
$ws.Range("N65").Value = -34881.428
$ws.Range("H107").Value = 882.43475
$ws.Range("I107").Value = 546.7646999999999
$ws.Range("K107").Value = 546.7646999999999
$ws.Range("M107").Value = 1373.2353
$ws.Range("H112").Value = 7555.3335
$ws.Range("J112").Value = 7555.3335
$ws.Range("L112").Value = 22666.0005
$ws.Range("N112").Value = -24882.0005
$ws.Range("H113").Value = 24365.6
$ws.Range("I113").Value = 29623.75
$ws.Range("K113").Value = 29623.75
$ws.Range("M113").Value = -26369.75
$ws.Range("H137").Value = 2083.9443
$ws.Range("I137").Value = 1756.8182
$ws.Range("J137").Value = 2598
$ws.Range("K137").Value = 5270.4546
$ws.Range("L137").Value = 7794
$ws.Range("M137").Value = -2720.4546
$ws.Range("N137").Value = -12894
$ws.Range("H138").Value = 6027.2363
$ws.Range("J138").Value = 6319.5
$ws.Range("L138").Value = 18958.5
$ws.Range("N138").Value = -29238.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3135.872
$ws.Range("I32").Value = 2561.8784
$ws.Range("J32").Value = 6675.5
$ws.Range("K32").Value = 2561.8784
$ws.Range("L32").Value = 6675.5
$ws.Range("M32").Value = -2274.8784
$ws.Range("N32").Value = -7249.5
$ws.Range("H45").Value = 1299.875
$ws.Range("I45").Value = 1223.2222
$ws.Range("J45").Value = 1398.4286
$ws.Range("K45").Value = 1223.2222
$ws.Range("L45").Value = 1398.4286
$ws.Range("M45").Value = -846.2221999999999
$ws.Range("N45").Value = -2152.4286
$ws.Range("H74").Value = 886.4583
$ws.Range("I74").Value = 808.7
$ws.Range("K74").Value = 808.7
$ws.Range("M74").Value = 65.29999999999995
$ws.Range("H77").Value = 886.4583
$ws.Range("I77").Value = 808.7
$ws.Range("K77").Value = 4043.5
$ws.Range("M77").Value = 324.5
$ws.Range("H132").Value = 2067.3103
$ws.Range("I132").Value = 1907.3334
$ws.Range("J132").Value = 2238.7144
$ws.Range("K132").Value = 5722.0002
$ws.Range("L132").Value = 6716.1432
$ws.Range("M132").Value = -3192.0002
$ws.Range("N132").Value = -11776.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1697.75
$ws.Range("I105").Value = 1741.75
$ws.Range("J105").Value = 1477.75
$ws.Range("K105").Value = 1741.75
$ws.Range("L105").Value = 1477.75
$ws.Range("M105").Value = 5.25
$ws.Range("N105").Value = -4971.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1321.9305
$ws.Range("J31").Value = 1535.25
$ws.Range("L31").Value = 1535.25
$ws.Range("N31").Value = -2125.25
$ws.Range("H34").Value = 1321.9305
$ws.Range("J34").Value = 1535.25
$ws.Range("L34").Value = 1535.25
$ws.Range("N34").Value = -1939.25
$ws.Range("H134").Value = 4372.7144
$ws.Range("I134").Value = 3722.2
$ws.Range("J134").Value = 5999
$ws.Range("K134").Value = 11166.6
$ws.Range("L134").Value = 17997
$ws.Range("M134").Value = -8631.599999999999
$ws.Range("N134").Value = -23067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H68").Value = 3446.9148
$ws.Range("I68").Value = 1399
$ws.Range("J68").Value = 4072.6667
$ws.Range("K68").Value = 4197
$ws.Range("L68").Value = 12218.0001
$ws.Range("M68").Value = -3386
$ws.Range("N68").Value = -13840.0001
$ws.Range("H71").Value = 3446.9148
$ws.Range("I71").Value = 1399
$ws.Range("J71").Value = 4072.6667
$ws.Range("K71").Value = 12591
$ws.Range("L71").Value = 36654.0003
$ws.Range("M71").Value = -8535
$ws.Range("N71").Value = -44766.0003
$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 12000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -10814
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 36000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -30072
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 4148.067
$ws.Range("J107").Value = 4230.0713
$ws.Range("L107").Value = 12690.2139
$ws.Range("N107").Value = -16530.2139
$ws.Range("H116").Value = 3000
$ws.Range("J116").Value = 3000
$ws.Range("L116").Value = 9000
$ws.Range("N116").Value = -15884
$ws.Range("H131").Value = 26354004
$ws.Range("I131").Value = 45454964
$ws.Range("J131").Value = 90182.125
$ws.Range("K131").Value = 136364892
$ws.Range("L131").Value = 270546.375
$ws.Range("M131").Value = -136359852
$ws.Range("N131").Value = -280626.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1134149.2
$ws.Range("I132").Value = 1604613.4
$ws.Range("J132").Value = 5035.6
$ws.Range("K132").Value = 4813840.199999999
$ws.Range("L132").Value = 15106.8
$ws.Range("M132").Value = -4811310.199999999
$ws.Range("N132").Value = -20166.8
$ws.Range("H141").Value = 73497.5
$ws.Range("J141").Value = 73497.5
$ws.Range("L141").Value = 73497.5
$ws.Range("N141").Value = -83857.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4373.227
$ws.Range("I136").Value = 3054.7036
$ws.Range("J136").Value = 6467.353
$ws.Range("K136").Value = 9164.110799999999
$ws.Range("L136").Value = 19402.059
$ws.Range("M136").Value = -6614.110799999999
$ws.Range("N136").Value = -24502.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1068.7576
$ws.Range("I132").Value = 821.9091
$ws.Range("J132").Value = 1562.4546
$ws.Range("K132").Value = 2465.7273
$ws.Range("L132").Value = 4687.3638
$ws.Range("M132").Value = 64.27269999999999
$ws.Range("N132").Value = -9747.363799999999
